# oferta.xlsx: refresh the promo subtext and let the "show" flag (B1)
# accept a plain 0/1 number instead of a =TRUE() boolean formula.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 ("subtext" row): reword/expand the order-window blurb.
$ws.Range("B3").Value = "Zamów dania na święta z odbiorem w barze. Przyjmujemy zamówienia do 15.04.2025 r., odbiór 18.04.2025 r. w godz. 9:30–18:00."

# B1 ("show" row): was `=TRUE()` rendered via a "TRUE"/"FALSE" custom number
# format; now a plain numeric literal (0 or 1) using the ordinary General
# number format, no formula.
$ws.Range("B1").Formula = "1"
$ws.Range("B1").NumberFormat = "General"
